$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '29.379.63'
$ws.Range('D2').Style = 'Normal'
$ws.Range('E2').Value = '  +0.14%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.884.97'
$ws.Range('D3').Style = 'Normal'
$ws.Range('E3').Value = '  +0.47%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.001'
$ws.Range('D4').Style = 'Normal'
$ws.Range('E4').Value = '  +0.10%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '0.7130'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  +0.18%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '242.47'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  +0.04%  '
$ws.Range('D7').NumberFormat = '@'
$ws.Range('D7').Value = '1.001'
$ws.Range('D7').Style = 'Normal'
$ws.Range('E7').Value = '  +0.07%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '0.08077'
$ws.Range('D8').Style = 'Normal'
$ws.Range('E8').Value = '  +4.14%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.3132'
$ws.Range('D9').Style = 'Normal'
$ws.Range('E9').Value = '  +0.66%  '
$ws.Range('E10').Value = '  +0.83%  '
$ws.Range('D11').NumberFormat = '@'
$ws.Range('D11').Value = '0.08359'
$ws.Range('D11').Style = 'Normal'
$ws.Range('E11').Value = '  -1.28%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '1.886.04'
$ws.Range('D12').Style = 'Normal'
$ws.Range('E12').Value = '  -0.25%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '0.7210'
$ws.Range('D13').Style = 'Normal'
$ws.Range('E13').Value = '  +1.35%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '5.246'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  +0.71%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '92.40'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  +1.11%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '6.277'
$ws.Range('D16').Style = 'Normal'
$ws.Range('E16').Value = '  +4.55%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '0.000008465'
$ws.Range('D17').Style = 'Normal'
$ws.Range('E17').Value = '  +2.04%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '29.381.67'
$ws.Range('D18').Style = 'Normal'
$ws.Range('E18').Value = '  +0.11%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '241.22'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -0.54%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '13.24'
$ws.Range('D20').Style = 'Normal'
$ws.Range('E20').Value = '  +0.14%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '2.126.34'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  +0.06%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '1.001'
$ws.Range('D22').Style = 'Normal'
$ws.Range('E22').Value = '  +0.11%  '
$ws.Range('D23').NumberFormat = '@'
$ws.Range('D23').Value = '7.827'
$ws.Range('D23').Style = 'Normal'
$ws.Range('E23').Value = '  -0.29%  '
$ws.Range('D24').NumberFormat = '@'
$ws.Range('D24').Value = '1.001'
$ws.Range('D24').Style = 'Normal'
$ws.Range('E24').Value = '  +0.05%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '0.1589'
$ws.Range('D25').Style = 'Normal'
$ws.Range('E25').Value = '  -1.36%  '
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '163.31'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  +0.43%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.074'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +0.59%  '
$ws.Range('E28').Value = '  +0.47%  '
$ws.Range('E29').Value = '  -0.40%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '4.425'
$ws.Range('D30').Style = 'Normal'
$ws.Range('E30').Value = '  +0.46%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '4.338'
$ws.Range('D31').Style = 'Normal'
$ws.Range('E31').Value = '  -0.16%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '1.205'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -5.52%  '
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '0.05379'
$ws.Range('D33').Style = 'Normal'
$ws.Range('E33').Value = '  +2.28%  '
$ws.Range('D34').NumberFormat = '@'
$ws.Range('D34').Value = '1.954'
$ws.Range('D34').Style = 'Normal'
$ws.Range('E34').Value = '  +1.10%  '
$ws.Range('E35').Value = '  +0.60%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '0.7503'
$ws.Range('D36').Style = 'Normal'
$ws.Range('E36').Value = '  +0.89%  '
$ws.Range('E37').Value = '  +0.62%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '0.01883'
$ws.Range('D38').Style = 'Normal'
$ws.Range('E38').Value = '  +0.79%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.285.21'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  +9.92%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '2.749'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +1.12%  '
$ws.Range('D41').NumberFormat = '@'
$ws.Range('D41').Value = '6.567'
$ws.Range('D41').Style = 'Normal'
$ws.Range('E41').Value = '  +2.82%  '
$ws.Range('D42').NumberFormat = '@'
$ws.Range('D42').Value = '73.53'
$ws.Range('D42').Style = 'Normal'
$ws.Range('E42').Value = '  +0.63%  '
$ws.Range('B43').Value = 'TrustWalletToken'
$ws.Range('C43').Value = 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt'
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '0.8930'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +0.60%  '
$ws.Range('B44').Value = 'Quant'
$ws.Range('C44').Value = 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt'
$ws.Range('D44').NumberFormat = '@'
$ws.Range('D44').Value = '110.43'
$ws.Range('D44').Style = 'Normal'
$ws.Range('E44').Value = '  +3.39%  '
$ws.Range('E45').Value = '  +0.08%  '
$ws.Range('E46').Value = '  +6.31%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '2.020.59'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -0.15%  '
$ws.Range('E48').Value = '  -0.48%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.5215'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  +0.21%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '9.502'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  +1.26%  '
$ws.Range('D51').NumberFormat = '@'
$ws.Range('D51').Value = '0.4392'
$ws.Range('D51').Style = 'Normal'
$ws.Range('E51').Value = '  +2.06%  '
